$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6451310861423221
$ws1.Range("C2").Value = 0.5848849945235487
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.7380787836903939
$ws1.Range("F2").Value = 0.875696949819613
$ws1.Range("G2").Value = 0.9734277501226951
$ws1.Range("H2").Value = 0.7119190899016679
$ws1.Range("I2").Value = 534
$ws1.Range("J2").Value = 379
$ws1.Range("K2").Value = 155
$ws1.Range("L2").Value = 0

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 ("0")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2902621722846442
$ws2.Range("D2").Value = 0.4499274310595066

# row 3 ("1")
$ws2.Range("B3").Value = 0.5848849945235487
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.7380787836903939

# row 4 ("accuracy")
$ws2.Range("B4").Value = 0.6451310861423221
$ws2.Range("C4").Value = 0.6451310861423221
$ws2.Range("D4").Value = 0.6451310861423221
$ws2.Range("E4").Value = 0.6451310861423221

# row 5 ("macro avg")
$ws2.Range("B5").Value = 0.7924424972617743
$ws2.Range("C5").Value = 0.6451310861423221
$ws2.Range("D5").Value = 0.5940031073749502

# row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.7924424972617743
$ws2.Range("C6").Value = 0.6451310861423221
$ws2.Range("D6").Value = 0.5940031073749503

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 ("Actual 0")
$ws3.Range("B2").Value = 155
$ws3.Range("C2").Value = 379

# row 3 ("Actual 1")
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 534
